$p = $ppt.ActivePresentation

# Remove the second slide (id 258) as per the diff
if ($p.Slides.Count -ge 2) {
    $p.Slides.Item(2).Delete()
}

$s = $p.Slides.Item(1)

# Shrink/reposition every shape on the slide to add a black margin for printing
$sh = $s.Shapes.Item(1)  # Picture 8
$sh.Left = 16.696537017822266
$sh.Top = 24.758582677165354
$sh.Width = 507.2183464566929
$sh.Height = 732.6369018554688

$sh = $s.Shapes.Item(2)  # Picture 11
$sh.Left = 231.07354330708662
$sh.Top = 109.20425196850394
$sh.Width = 74.84055118110236
$sh.Height = 74.84055118110236

$sh = $s.Shapes.Item(3)  # Rounded Rectangle 9
$sh.Left = 49.253385826771655
$sh.Top = 228.9703937007874
$sh.Width = 435.9438781738281
$sh.Height = 500.37724409448816

$sh = $s.Shapes.Item(4)  # TextBox 3
$sh.Left = 59.523231506347656
$sh.Top = 43.58338928222656
$sh.Width = 425.67393700787403
$sh.Height = 70.37039947509766

$sh = $s.Shapes.Item(5)  # TextBox 5
$sh.Left = 58.34007874015748
$sh.Top = 166.08567810058594
$sh.Width = 462.321044921875
$sh.Height = 52.21031496062992

$sh = $s.Shapes.Item(6)  # TextBox 4
$sh.Left = 72.90158081054688
$sh.Top = 115.96629921259843
$sh.Width = 398.9172668457031
$sh.Height = 38.59023622047244

$sh = $s.Shapes.Item(7)  # TextBox 10
$sh.Left = 64.10488891601562
$sh.Top = 286.8895568847656
$sh.Width = 422.5654602050781
$sh.Height = 415.4123840332031

$sh = $s.Shapes.Item(8)  # Picture 13
$sh.Left = 15.230000495910645
$sh.Top = 208.15456692913386
$sh.Width = 507.21844482421875
$sh.Height = 542.8248818897638

$sh = $s.Shapes.Item(9)  # TextBox 12
$sh.Left = 50.774015748031495
$sh.Top = 226.03818897637797
$sh.Width = 126.8939437866211
$sh.Height = 38.59023622047244
